$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B ("R40") is being replaced with the text value "1".
# Force the cell to remain text (not get auto-converted to a number) by
# setting the number format to Text before assigning the new value.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
